$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1939546599496222
$ws.Range("C2").Value = 0.5642317380352645
$ws.Range("J2").Value = 0.01763224181360202
$ws.Range("P2").Value = 0.1360201511335013
$ws.Range("S2").Value = 0.08816120906801007
$ws.Range("B3").Value = 0.01310043668122271
$ws.Range("C3").Value = 0.01310043668122271
$ws.Range("J3").Value = 0.03930131004366812
$ws.Range("P3").Value = 0.7336244541484717
$ws.Range("S3").Value = 0.2008733624454148
$ws.Range("J4").Value = 0.08064516129032258
$ws.Range("P4").Value = 0.6451612903225806
$ws.Range("S4").Value = 0.2741935483870968
$ws.Range("B6").Value = 0.03888888888888889
$ws.Range("D6").Value = 0.01666666666666667
$ws.Range("E6").Value = 0.005555555555555556
$ws.Range("F6").Value = 0.04444444444444445
$ws.Range("J6").Value = 0.3333333333333333
$ws.Range("O6").Value = 0.005555555555555556
$ws.Range("Q6").Value = 0.1222222222222222
$ws.Range("R6").Value = 0.03888888888888889
$ws.Range("S6").Value = 0.3944444444444444
$ws.Range("B7").Value = 0.1245283018867925
$ws.Range("D7").Value = 0.01509433962264151
$ws.Range("F7").Value = 0.03773584905660377
$ws.Range("J7").Value = 0.1320754716981132
$ws.Range("O7").Value = 0.01509433962264151
$ws.Range("Q7").Value = 0.2188679245283019
$ws.Range("R7").Value = 0.06037735849056604
$ws.Range("S7").Value = 0.3962264150943396
$ws.Range("B8").Value = 0.1653944020356234
$ws.Range("D8").Value = 0.02290076335877863
$ws.Range("F8").Value = 0.06106870229007633
$ws.Range("J8").Value = 0.1246819338422392
$ws.Range("O8").Value = 0.01526717557251908
$ws.Range("Q8").Value = 0.1577608142493639
$ws.Range("R8").Value = 0.04580152671755725
$ws.Range("S8").Value = 0.4071246819338422
$ws.Range("B9").Value = 0.155
$ws.Range("D9").Value = 0.035
$ws.Range("F9").Value = 0.08
$ws.Range("J9").Value = 0.095
$ws.Range("O9").Value = 0.02
$ws.Range("Q9").Value = 0.145
$ws.Range("R9").Value = 0.06
$ws.Range("S9").Value = 0.41
$ws.Range("B10").Value = 0.1404893449092344
$ws.Range("D10").Value = 0.03078137332280979
$ws.Range("E10").Value = 0.001578531965272297
$ws.Range("F10").Value = 0.05209155485398579
$ws.Range("J10").Value = 0.1207576953433307
$ws.Range("O10").Value = 0.01499605367008682
$ws.Range("Q10").Value = 0.1949486977111287
$ws.Range("R10").Value = 0.03867403314917127
$ws.Range("S10").Value = 0.4056827150749803
$ws.Range("G11").Value = 0.1358024691358025
$ws.Range("J11").Value = 0.0691358024691358
$ws.Range("K11").Value = 0.182716049382716
$ws.Range("L11").Value = 0.5950617283950618
$ws.Range("S11").Value = 0.01728395061728395
$ws.Range("G12").Value = 0.7459677419354839
$ws.Range("J12").Value = 0.2056451612903226
$ws.Range("K12").Value = 0.004032258064516129
$ws.Range("L12").Value = 0.01612903225806452
$ws.Range("S12").Value = 0.0282258064516129
$ws.Range("G13").Value = 0.6739130434782609
$ws.Range("J13").Value = 0.2173913043478261
$ws.Range("S13").Value = 0.108695652173913
$ws.Range("F15").Value = 0.004807692307692308
$ws.Range("H15").Value = 0.1490384615384615
$ws.Range("I15").Value = 0.08653846153846154
$ws.Range("J15").Value = 0.3461538461538461
$ws.Range("K15").Value = 0.07211538461538461
$ws.Range("M15").Value = 0.02403846153846154
$ws.Range("N15").Value = 0.004807692307692308
$ws.Range("O15").Value = 0.0625
$ws.Range("S15").Value = 0.25
$ws.Range("F16").Value = 0.01185770750988142
$ws.Range("H16").Value = 0.1541501976284585
$ws.Range("I16").Value = 0.06324110671936758
$ws.Range("J16").Value = 0.383399209486166
$ws.Range("K16").Value = 0.1620553359683795
$ws.Range("M16").Value = 0.03162055335968379
$ws.Range("O16").Value = 0.08300395256916997
$ws.Range("S16").Value = 0.1106719367588933
$ws.Range("F17").Value = 0.01927710843373494
$ws.Range("H17").Value = 0.1373493975903614
$ws.Range("I17").Value = 0.09156626506024096
$ws.Range("J17").Value = 0.4144578313253012
$ws.Range("K17").Value = 0.1301204819277108
$ws.Range("M17").Value = 0.01686746987951807
$ws.Range("O17").Value = 0.05783132530120482
$ws.Range("S17").Value = 0.1325301204819277
$ws.Range("F18").Value = 0.05882352941176471
$ws.Range("H18").Value = 0.1470588235294118
$ws.Range("I18").Value = 0.07843137254901961
$ws.Range("J18").Value = 0.2941176470588235
$ws.Range("K18").Value = 0.1274509803921569
$ws.Range("M18").Value = 0.02941176470588235
$ws.Range("O18").Value = 0.07843137254901961
$ws.Range("S18").Value = 0.1862745098039216
$ws.Range("F19").Value = 0.01392961876832845
$ws.Range("H19").Value = 0.187683284457478
$ws.Range("I19").Value = 0.08870967741935484
$ws.Range("J19").Value = 0.3563049853372434
$ws.Range("K19").Value = 0.1495601173020528
$ws.Range("M19").Value = 0.01832844574780059
$ws.Range("S19").Value = 0.1282991202346041
